$wb = $excel.ActiveWorkbook

# Citywide Totals (sheet1)
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 4751
$ws.Range("J3").Value = 8079
$ws.Range("K3").Value = 4887
$ws.Range("C4").Value = 1850
$ws.Range("H4").Value = 1738
$ws.Range("K4").Value = 1009
$ws.Range("K5").Value = 350
$ws.Range("K6").Value = 5500
$ws.Range("C7").Value = 28395
$ws.Range("H7").Value = 26051
$ws.Range("J7").Value = 29299
$ws.Range("K7").Value = 16497

# By Neighborhood (sheet2)
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K7").Value = 484
$ws.Range("K8").Value = 1111
$ws.Range("K11").Value = 322
$ws.Range("K12").Value = 31
$ws.Range("K14").Value = 91
$ws.Range("K18").Value = 111
$ws.Range("K19").Value = 487
$ws.Range("K20").Value = 379
$ws.Range("K23").Value = 168
$ws.Range("K25").Value = 80
$ws.Range("K27").Value = 150
$ws.Range("K29").Value = 877
$ws.Range("K30").Value = 61
$ws.Range("K33").Value = 695
$ws.Range("K34").Value = 87
$ws.Range("K36").Value = 213
$ws.Range("J37").Value = 902
$ws.Range("K37").Value = 555
$ws.Range("K42").Value = 613
$ws.Range("K43").Value = 145
$ws.Range("K47").Value = 109
$ws.Range("K48").Value = 207
$ws.Range("K51").Value = 207
$ws.Range("K52").Value = 432
$ws.Range("K54").Value = 324
$ws.Range("K55").Value = 188
$ws.Range("K57").Value = 57
$ws.Range("C63").Value = 279
$ws.Range("H63").Value = 289
$ws.Range("K63").Value = 55
$ws.Range("K67").Value = 636
$ws.Range("K68").Value = 43
$ws.Range("K78").Value = 194
$ws.Range("K83").Value = 358
$ws.Range("K85").Value = 751
$ws.Range("K86").Value = 110
$ws.Range("K89").Value = 232
$ws.Range("K90").Value = 151
$ws.Range("K93").Value = 61
$ws.Range("K94").Value = 216
$ws.Range("K95").Value = 291
$ws.Range("K96").Value = 178
$ws.Range("K98").Value = 82
$ws.Range("K99").Value = 276
$ws.Range("C101").Value = 28395
$ws.Range("H101").Value = 26051
$ws.Range("J101").Value = 29299
$ws.Range("K101").Value = 16497

# Bridgeport (sheet3)
$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("K6").Value = 33
$ws.Range("K7").Value = 91

# West Ridge (sheet4)
$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K2").Value = 57
$ws.Range("K7").Value = 178

# Auburn Gresham (sheet5)
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 168
$ws.Range("K6").Value = 125
$ws.Range("K7").Value = 484

# Belmont Cragin (sheet6)
$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K3").Value = 85
$ws.Range("K7").Value = 322

# Uptown (sheet7)
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K4").Value = 26
$ws.Range("K7").Value = 232

# South Shore (sheet8)
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K3").Value = 253
$ws.Range("K4").Value = 44
$ws.Range("K6").Value = 174
$ws.Range("K7").Value = 751

# Little Village (sheet9)
$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K3").Value = 116
$ws.Range("K6").Value = 161
$ws.Range("K7").Value = 432

# Austin (sheet12)
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 310
$ws.Range("K3").Value = 334
$ws.Range("K6").Value = 373
$ws.Range("K7").Value = 1111

# South Chicago (sheet13)
$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K3").Value = 128
$ws.Range("K4").Value = 19
$ws.Range("K6").Value = 81
$ws.Range("K7").Value = 358

# Garfield Park (sheet14)
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K3").Value = 262
$ws.Range("K7").Value = 695

# West Pullman (sheet15)
$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K3").Value = 102
$ws.Range("K7").Value = 291

# Grand Crossing (sheet16)
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J3").Value = 302
$ws.Range("K3").Value = 182
$ws.Range("K4").Value = 25
$ws.Range("K6").Value = 165
$ws.Range("J7").Value = 902
$ws.Range("K7").Value = 555

# New City (sheet17)
$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 115
$ws.Range("K6").Value = 145

# Woodlawn (sheet18)
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K2").Value = 71
$ws.Range("K3").Value = 113
$ws.Range("K7").Value = 276

# Fuller Park (sheet19)
$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("K3").Value = 22
$ws.Range("K7").Value = 61

# North Lawndale (sheet21)
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 183
$ws.Range("K3").Value = 220
$ws.Range("K5").Value = 13
$ws.Range("K7").Value = 636

# Loop (sheet24)
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K6").Value = 172
$ws.Range("K7").Value = 324

# Englewood (sheet25)
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 250
$ws.Range("K3").Value = 315
$ws.Range("K7").Value = 877

# Lake View (sheet26)
$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K2").Value = 28
$ws.Range("K6").Value = 102
$ws.Range("K7").Value = 207

# Chatham (sheet27)
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 146
$ws.Range("K6").Value = 154
$ws.Range("K7").Value = 487

# Humboldt Park (sheet32)
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K3").Value = 191
$ws.Range("K6").Value = 230
$ws.Range("K7").Value = 613

# Rogers Park (sheet35)
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K2").Value = 57
$ws.Range("K4").Value = 21
$ws.Range("K7").Value = 194

# Lower West Side (sheet36)
$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K3").Value = 55
$ws.Range("K7").Value = 188

# Douglas (sheet39)
$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K2").Value = 46
$ws.Range("K7").Value = 168

# Chicago Lawn (sheet44)
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 128
$ws.Range("K7").Value = 379

# Calumet Heights (sheet45)
$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K6").Value = 32
$ws.Range("K7").Value = 111

# Grand Boulevard (sheet47)
$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K6").Value = 51
$ws.Range("K7").Value = 213

# West Lawn (sheet48)
$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 61

# Garfield Ridge (sheet50)
$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K2").Value = 29
$ws.Range("K6").Value = 28
$ws.Range("K7").Value = 87

# West Loop (sheet51)
$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K2").Value = 65
$ws.Range("K6").Value = 91
$ws.Range("K7").Value = 216

# East Side (sheet52)
$ws = $wb.Worksheets.Item("East Side")
$ws.Range("K6").Value = 17
$ws.Range("K7").Value = 80

# Kenwood (sheet53)
$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K3").Value = 33
$ws.Range("K6").Value = 37
$ws.Range("K7").Value = 109

# Wicker Park (sheet55)
$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K2").Value = 12
$ws.Range("K7").Value = 82

# Edgewater (sheet71)
$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K3").Value = 36
$ws.Range("K6").Value = 57
$ws.Range("K7").Value = 150

# Streeterville (sheet72)
$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("K4").Value = 47
$ws.Range("K7").Value = 110

# Washington Heights (sheet74)
$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K3").Value = 48
$ws.Range("K7").Value = 151

# Little Italy, UIC (sheet75)
$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K2").Value = 59
$ws.Range("K7").Value = 207

# North Park (sheet76)
$ws = $wb.Worksheets.Item("North Park")
$ws.Range("K2").Value = 18
$ws.Range("K7").Value = 43

# Mckinley Park (sheet77)
$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("K2").Value = 17
$ws.Range("K4").Value = 3
$ws.Range("K7").Value = 57

# Hyde Park (sheet79)
$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K3").Value = 39
$ws.Range("K7").Value = 145

# Beverly (sheet91)
$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("K6").Value = 9
$ws.Range("K7").Value = 31
